$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 102.1070184418289
$ws.Range("O2").Value = 0.0011
$ws.Range("Q2").Value = 0.0008
$ws.Range("S2").Value = 0.0008
$ws.Range("T2").Value = 0.0016
$ws.Range("O3").Value = 0.001
$ws.Range("P3").Value = 0.0011
$ws.Range("Q3").Value = 0.001
$ws.Range("R3").Value = 0.001
$ws.Range("O4").Value = 0.0013
$ws.Range("P4").Value = 0.0012
$ws.Range("R4").Value = 0.0011
$ws.Range("S4").Value = 0.001
$ws.Range("I5").Value = 39.37213886457447
$ws.Range("J5").Value = 45.69297498377761
$ws.Range("K5").Value = 49.70847756767837
$ws.Range("M5").Value = 31.22691510942754
$ws.Range("Q5").Value = 0.0012
$ws.Range("S5").Value = 0.001
$ws.Range("I6").Value = 37.00463120931033
$ws.Range("J6").Value = 44.27627586887465
$ws.Range("K6").Value = 41.01081695451413
$ws.Range("L6").Value = 31.88252949105588
$ws.Range("M6").Value = 32.11184821924887
$ws.Range("N6").Value = 31.45623383762054
$ws.Range("O6").Value = 0.0021
$ws.Range("P6").Value = 0.0023
$ws.Range("R6").Value = 0.0018
$ws.Range("S6").Value = 0.0017
$ws.Range("T6").Value = 0.0016
$ws.Range("I7").Value = 42.19843329089865
$ws.Range("J7").Value = 40.31098556893893
$ws.Range("K7").Value = 42.57832116775565
$ws.Range("L7").Value = 31.88252949105588
$ws.Range("O7").Value = 0.0041
$ws.Range("I8").Value = 560654.4810176021
$ws.Range("J8").Value = 572266.7876941945
$ws.Range("K8").Value = 565783.7988525628
$ws.Range("M8").Value = 49294.74163904427
$ws.Range("O8").Value = 0.0085
$ws.Range("P8").Value = 0.0084
$ws.Range("Q8").Value = 0.0086
$ws.Range("R8").Value = 0.0688
$ws.Range("S8").Value = 0.0685
$ws.Range("T8").Value = 0.0683
$ws.Range("I9").Value = 535498.8765180835
$ws.Range("J9").Value = 535156.0096518323
$ws.Range("K9").Value = 530950.212590666
$ws.Range("L9").Value = 48359.58007325514
$ws.Range("M9").Value = 48350.55885416948
$ws.Range("N9").Value = 48554.11488899639
$ws.Range("O9").Value = 0.0347
$ws.Range("P9").Value = 0.0357
$ws.Range("Q9").Value = 0.0345
$ws.Range("R9").Value = 0.331
$ws.Range("S9").Value = 0.3248
$ws.Range("T9").Value = 0.3297
$ws.Range("I10").Value = 513333.8967455443
$ws.Range("J10").Value = 535586.108007305
$ws.Range("K10").Value = 527628.7018299692
$ws.Range("L10").Value = 48928.00065664363
$ws.Range("N10").Value = 49143.7729793856
$ws.Range("O10").Value = 0.068
$ws.Range("P10").Value = 0.0696
$ws.Range("Q10").Value = 0.0678
$ws.Range("R10").Value = 0.6638
$ws.Range("S10").Value = 0.6497
$ws.Range("T10").Value = 0.6545
